# Update the cached "datetimeFigureOut" date field text (footer placeholder)
# from 8/30/2021 to 9/28/2021 across the slide master and every slide layout,
# and bump the version string shown on the splash slide from "version 1.6"
# to "version 1.6.1".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                if ($shape.TextFrame.TextRange.Text -eq "8/30/2021") {
                    $shape.TextFrame.TextRange.Text = "9/28/2021"
                }
            }
        }
    }
}

# Slide master's Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout's Date Placeholder.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DatePlaceholder $layout.Shapes
}

# Bump the version number shown on the splash slide.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $textRange = $shape.TextFrame.TextRange
        $paraCount = $textRange.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            # Each non-final paragraph's .Text carries a trailing "\r"
            # paragraph-mark, so compare against the trimmed value.
            $para = $textRange.Paragraphs($j)
            if ($para.Text.TrimEnd("`r") -eq "version 1.6") {
                # Stage through an unrelated placeholder string first so the
                # final assignment shares no common prefix with the old run
                # text; that keeps the run formatting intact as a single
                # <a:r> (instead of being split into an unchanged-prefix run
                # plus a new ".1" run) while still only touching this run.
                $para.Text = "__TMP_VERSION_PLACEHOLDER__"
                $para = $textRange.Paragraphs($j)
                $para.Text = "version 1.6.1"
            }
        }
    }
}
